$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.556.47'
$ws.Range("E2").Value = '  +2.51%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.184.45'
$ws.Range("E3").Value = '  +0.08%  '

$ws.Range("E4").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.39'
$ws.Range("E5").Value = '  +1.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.24'
$ws.Range("E6").Value = '  +10.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  +1.26%  '

$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.590'
$ws.Range("E9").Value = '  +0.65%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.95'
$ws.Range("E10").Value = '  +6.24%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0916'
$ws.Range("E11").Value = '  -0.08%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.93'
$ws.Range("E12").Value = '  +2.29%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.103'
$ws.Range("E13").Value = '  +1.92%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.509.61'
$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.21'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.170.96'
$ws.Range("E16").Value = '  -0.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.771'
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.463.53'
$ws.Range("E18").Value = '  +2.49%  '

$ws.Range("E19").Value = '  +0.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.81'
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.89'
$ws.Range("E21").Value = '  +0.23%  '

$ws.Range("E22").Value = '  +11.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.80'
$ws.Range("E23").Value = '  +1.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.80'
$ws.Range("E24").Value = '  -5.91%  '

$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '42.48'
$ws.Range("E26").Value = '  +15.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.62'
$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("E28").Value = '  -0.43%  '

$ws.Range("E29").Value = '  +0.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.22'
$ws.Range("E30").Value = '  +1.52%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '172.93'
$ws.Range("E31").Value = '  +1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.33'
$ws.Range("E32").Value = '  +1.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0868'
$ws.Range("E33").Value = '  +7.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.27'
$ws.Range("E34").Value = '  +2.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.113'
$ws.Range("E35").Value = '  +5.25%  '

$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.44'
$ws.Range("E37").Value = '  +4.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0350'
$ws.Range("E38").Value = '  +3.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.93'
$ws.Range("E39").Value = '  +9.38%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.82'
$ws.Range("E40").Value = '  +13.99%  '

$ws.Range("E41").Value = '  +1.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '62.51'
$ws.Range("E42").Value = '  +5.42%  '

$ws.Range("E43").Value = '  +5.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.198'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '100.49'
$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0978'
$ws.Range("E46").Value = '  +0.93%  '

$ws.Range("E47").Value = '  -0.83%  '

$ws.Range("E48").Value = '  +3.82%  '

$ws.Range("E49").Value = '  +1.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.53'
$ws.Range("E50").Value = '  +26.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.435'
$ws.Range("E51").Value = '  -6.88%  '
